$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for a new "Choice" row below the current last data row (9),
#    which currently carries the thick bottom border / bold look because it
#    is the last row of the table. Insert a new row at 10 (pushes the old
#    trailing blank row down to 11).
$ws.Rows.Item(10).Insert()

# 3. Copy row 9's current (bottom-of-table) formatting down into the new
#    row 10 so row 10 becomes the new visual bottom of the table.
$ws.Range("C9:F9").Copy() | Out-Null
$ws.Range("C10:F10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# 4. Restore row 9 back to a plain interior row, matching rows 6-8.
$ws.Range("C8:F8").Copy() | Out-Null
$ws.Range("C9:F9").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws.Application.CutCopyMode = $false

# 5. Fill in row 9 values/formulas (unchanged data, just re-set after the
#    format copy wiped the cell contents).
$ws.Cells.Item(9, 3).Value = "Narrow rule (Logit)"
$ws.Cells.Item(9, 4).Formula = "=ROUND(H9,2)"
$ws.Cells.Item(9, 5).Formula = "=ROUND(I9,2)"
$ws.Cells.Item(9, 6).Formula = "=SUM(D9:E9)"
$ws.Cells.Item(9, 8).Value = 9.3628444671630859
$ws.Cells.Item(9, 9).Value = 10.242380142211911

# 6. Fill in the new "Choice" row 10.
$ws.Cells.Item(10, 3).Value = "Choice"
$ws.Cells.Item(10, 4).Formula = "=ROUND(H10,2)"
$ws.Cells.Item(10, 5).Formula = "=ROUND(I10,2)"
$ws.Cells.Item(10, 6).Formula = "=ROUND(J10,2)"
$ws.Cells.Item(10, 8).Value = 80.296737670898438
$ws.Cells.Item(10, 9).Value = 21.367921829223629
$ws.Cells.Item(10, 10).Value = 71.23077392578125

$ws.Range("C3:J10").Calculate() | Out-Null
